# Add a new "Slovakia" worksheet (cloned from "Portugal") with its own market data.
$wb = $excel.ActiveWorkbook

$portugal = $wb.Worksheets.Item("Portugal")

# Copy Portugal to the end of the workbook -> new sheet becomes active.
$portugal.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))

$slovakia = $wb.Worksheets.Item($wb.Worksheets.Count)
$slovakia.Name = "Slovakia"

# Update the market-specific values on the new sheet.
$slovakia.Range("B2").Value = "Slovakia Market"
$slovakia.Range("B4").Value = "NGC-2930/T3222"

# Rows 3 & 4 shrink back to the default height (only row 5 keeps the taller,
# wrapped "Expected value" row), unlike the row heights inherited from Portugal.
$slovakia.Rows.Item(3).AutoFit()
$slovakia.Rows.Item(4).AutoFit()

# Portugal is no longer the selected tab; its cursor moves to B12.
$portugal.Range("B12").Select()

# Selection / active cell on the new sheet, which remains the active tab.
$slovakia.Range("A11").Select()
